# The workbook tracks weekly "Puerro" (leek) price records for Vega Modelo
# de Temuco, newest record first. A new weekly record was added at row 41,
# pushing all the existing records (old rows 41-177) down by one row
# (to new rows 42-178) while keeping rows 2-40 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 41; Excel automatically shifts rows 41:177
# down to 42:178 (values, formulas and styles all move together), and
# extends the sheet's used range/dimension accordingly.
$ws.Rows("41:41").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44607
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112005
$ws.Range("G41").Value = "Puerro"
$ws.Range("H41").Value = "Azul de Maquehue"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 30
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = 10000
$ws.Range("N41").Value = "$/docena de paquetes"
$ws.Range("O41").Value = "Provincia de Cautín"
$ws.Range("P41").Value = 833
$ws.Range("Q41").Value = 12
$ws.Range("R41").Value = "Hortaliza"
